$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "gnjdF268"
$ws.Range("C2").Value = "eospvak97"
$ws.Range("D2").Value = "mp7!5PU&"
$ws.Range("F2").Value = "hrBndtCO"
$ws.Range("G2").Value = "yWoz"
$ws.Range("B2").Value = 231006170
